# Update zhongshu_wangge.xlsx for 2021 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 创业板50（159949） block (rows 7-9): refresh 中枢0 / 中枢1 numbers ---
$ws.Range("B8").Value = "1.169/1.201"
$ws.Range("D8").Value = "1.114/1.150"
$ws.Range("B9").Value = "1.221/1.250"
$ws.Range("D9").Value = "1.176/1.199"

# --- 300ETF（510300） block (rows 16-18): refresh numbers, clear stale execution flags ---
$ws.Range("D17").Value = "5.080/5.121"
$ws.Range("E17").ClearContents()
$ws.Range("B18").Value = "5.112/5.152"
$ws.Range("C18").ClearContents()
$ws.Range("D18").Value = "5.043/5.097"

# --- 科创50（588000） block (rows 25-27): refresh numbers, clear stale execution flags ---
$ws.Range("C26").ClearContents()
$ws.Range("D26").ClearContents()
$ws.Range("B27").Value = "1.408/1.420"
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("G27").ClearContents()

# --- 证券ETF（512880） block (rows 34-36): clear stale execution flags ---
$ws.Range("C35").ClearContents()
$ws.Range("E35").ClearContents()
$ws.Range("C36").ClearContents()
$ws.Range("E36").ClearContents()

# --- Update the saved view/selection state ---
$ws.Range("E38:E39").Select()
